$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "customer_name"
$ws.Range("B9").Value = "No"
$ws.Range("C9").Value = "string"
$ws.Range("K9").Value = "Test User"

$ws.Range("A10").Value = "customer_email"
$ws.Range("B10").Value = "Yes"
$ws.Range("C10").Value = "string"
$ws.Range("K10").Value = "test@test.com"

$ws.Range("A4:C4").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)
$ws.Range("K4").Copy()
$ws.Range("K9").PasteSpecial(-4122)

$ws.Range("A4:C4").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122)
$ws.Range("K4").Copy()
$ws.Range("K10").PasteSpecial(-4122)
